# Updated symbol list on Tue Dec 27 22:45:59 UTC 2022 with GitHub Actions
# Applies the latest coinranking.com price/volume snapshot to Sheet1.
#
# Price values in column D are stored as text (they include values like
# "--" for unavailable prices), so every write is prefixed with a leading
# apostrophe to force Excel to keep them as literal text instead of
# auto-converting to a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

# Column D: Price
Set-TextValue "D2"  "246.08"
Set-TextValue "D3"  "23.97"
Set-TextValue "D4"  "5.354"
Set-TextValue "D5"  "0.05812"
Set-TextValue "D6"  "3.375"
Set-TextValue "D7"  "6.471"
Set-TextValue "D8"  "0.8094"
Set-TextValue "D9"  "0.9202"
Set-TextValue "D11" "0.07387"
Set-TextValue "D12" "0.03193"
Set-TextValue "D14" "0.09378"
Set-TextValue "D15" "3.864"
Set-TextValue "D16" "0.001567"
Set-TextValue "D17" "0.04697"
Set-TextValue "D18" "0.0005996"
Set-TextValue "D19" "0.006009"
Set-TextValue "D21" "0.004703"
Set-TextValue "D28" "0.0002348"
Set-TextValue "D41" "0.006365"
Set-TextValue "D43" "0.002899"
Set-TextValue "D44" "0.009051"
Set-TextValue "D45" "0.00005241"
Set-TextValue "D47" "0.6850"

# Column E: Volume(1h) label (rank + name + symbol [+ Best/Worst in 24h tag])
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

Write-Host "Applied symbol-list update."
